$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the daily spot-price record (new day + refreshed hourly prices)
$ws.Range("A2").Value = 45909

$ws.Range("B2").Value = 63.33
$ws.Range("C2").Value = 63.33
$ws.Range("D2").Value = 61.51
$ws.Range("E2").Value = 61.51
$ws.Range("F2").Value = 61.51
$ws.Range("G2").Value = 69.45
$ws.Range("H2").Value = 83.09
$ws.Range("I2").Value = 96.09999999999999
$ws.Range("J2").Value = 96.44
$ws.Range("K2").Value = 80.70999999999999
$ws.Range("L2").Value = 54.93
$ws.Range("M2").Value = 30.22
$ws.Range("N2").Value = 26
$ws.Range("O2").Value = 19.25
$ws.Range("P2").Value = 17.19
$ws.Range("Q2").Value = 8.91
$ws.Range("R2").Value = 16.57
$ws.Range("S2").Value = 30
$ws.Range("T2").Value = 56
$ws.Range("U2").Value = 90
$ws.Range("V2").Value = 111.37
$ws.Range("W2").Value = 120
$ws.Range("X2").Value = 110
$ws.Range("Y2").Value = 98.06999999999999
$ws.Range("Z2").Value = 63.56

$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 109.86
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 115.68
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 104.04
$ws.Range("AG2").Value = "0h-18h"
